$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for rows 2-23 (columns A-F), reflecting the reordered
# (and one corrected) data rows from calcWaysWinsInReels2.
$data = @(
    @(901, 16, 15, 45, 60, 60),
    @(301, 6, 45, 30, 60, 45),
    @(1201, 2, 10, 10, 10, 10),
    @(1202, 2, 10, 10, 10, 10),
    @(101, 9, 30, 15, 60, 15),
    @(401, 9, 48, 67, 75, 45),
    @(701, 3, 90, 45, 97, 15),
    @(601, 9, 60, 67, 60, 42),
    @(801, 3, 67, 65, 52, 45),
    @(501, 9, 52, 30, 75, 45),
    @(1203, 3, 15, 15, 15, 15),
    @(902, 1, 0, 0, 0, 0),
    @(1001, 18, 30, 75, 60, 72),
    @(201, 9, 30, 15, 45, 30),
    @(802, 0, 4, 5, 4, 0),
    @(502, 0, 4, 0, 0, 0),
    @(1101, 0, 15, 30, 30, 0),
    @(3, 0, 3, 3, 3, 3),
    @(1, 0, 2, 2, 2, 2),
    @(2, 0, 2, 2, 2, 2),
    @(602, 0, 0, 4, 0, 9),
    @(402, 0, 0, 4, 0, 0)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}

$wb.Save()
